# login() is added with negative testcases
# This script rebuilds the "Login" worksheet to match the target layout:
# header rows, new negative-scenario rows, a new bold+blue sub-header style,
# and a set of mailto: hyperlinks on the email-looking cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Login")

# ---- Row 1-2: headers (unchanged text, A1 text changes) ----
$ws.Cells.Item(1,1).Value = "LoginButton disability"

$ws.Cells.Item(2,1).Value = "Login ID"
$ws.Cells.Item(2,2).Value = "Password"

# ---- Row 3: both fields blank (two spaces) ----
$ws.Cells.Item(3,1).Value = "  "
$ws.Cells.Item(3,2).Value = "  "
$ws.Cells.Item(3,3).Value = "two blankspaces in each cell (both the fields are empty)"

# ---- Row 4: only email present ----
$ws.Cells.Item(4,1).Value = "test3456"
$ws.Cells.Item(4,2).Value = "styletag123"

# ---- Row 5 ----
$ws.Cells.Item(5,1).Value = "test3456@"
$ws.Cells.Item(5,2).Value = "styletag123"

# ---- Row 6 (hyperlinked partial email) ----
$ws.Cells.Item(6,1).Value = "test3456@st"
$ws.Cells.Item(6,2).Value = "styletag123"

# ---- Row 7 (hyperlinked partial email) ----
$ws.Cells.Item(7,1).Value = "test3456@styletag."
$ws.Cells.Item(7,2).Value = "styletag123"

# ---- Row 8 (hyperlinked partial email) ----
$ws.Cells.Item(8,1).Value = "test3456@styletag.c"
$ws.Cells.Item(8,2).Value = "styletag123"

# ---- Row 9: blank email (6 spaces), blank note ----
$ws.Cells.Item(9,1).Value = "      "
$ws.Cells.Item(9,2).Value = "styletag123"
$ws.Cells.Item(9,3).Value = "blank emailis"

# ---- Row 10: full valid email, blank password ----
$ws.Cells.Item(10,1).Value = "test3456@styletag.com"
$ws.Cells.Item(10,2).Value = "        "
$ws.Cells.Item(10,3).Value = "blank passwors"

# ---- Row 11 ----
$ws.Cells.Item(11,1).Value = "Test3456@.com"
$ws.Cells.Item(11,2).Value = "styletag123"

# ---- Row 12 ----
$ws.Cells.Item(12,1).Value = "@.com"
$ws.Cells.Item(12,2).Value = "styletag123"

# ---- Row 13: new bold+blue sub header ----
$ws.Cells.Item(13,1).Value = "LoginButton enabled"
$ws.Cells.Item(13,1).Font.Bold = $true
$ws.Cells.Item(13,1).Font.Color = 16711680
$ws.Cells.Item(13,2).Value = "following data are not valid for login"

# ---- Row 14: valid email, invalid password ----
$ws.Cells.Item(14,1).Value = "test3456@styletag.com"
$ws.Cells.Item(14,2).Value = "styletag12"
$ws.Cells.Item(14,3).Value = "invalid password"

# ---- Row 15: not registered credentials ----
$ws.Cells.Item(15,1).Value = "test34567@styletag.com"
$ws.Cells.Item(15,2).Value = "styletag123"
$ws.Cells.Item(15,3).Value = "not registered credentials"

# ---- Row 16: Valid data sub header ----
$ws.Cells.Item(16,1).Value = "Valid data"

# ---- Row 17: registered data ----
$ws.Cells.Item(17,1).Value = "test3456@styletag.com"
$ws.Cells.Item(17,2).Value = "styletag123"
$ws.Cells.Item(17,3).Value = "Registered data"

# ---- Hyperlinks on the email-looking cells ----
$ws.Hyperlinks.Add($ws.Cells.Item(6,1), "mailto:test3456@styletag.com", "", "", "test3456@st")
$ws.Cells.Item(6,1).Font.Underline = $false
$ws.Cells.Item(6,1).Font.Color = 16711680
$ws.Cells.Item(6,1).Font.Name = "Arial"

$ws.Hyperlinks.Add($ws.Cells.Item(7,1), "mailto:test3456@styletag.com", "", "", "test3456@styletag.")
$ws.Cells.Item(7,1).Font.Underline = $false
$ws.Cells.Item(7,1).Font.Color = 16711680
$ws.Cells.Item(7,1).Font.Name = "Arial"

$ws.Hyperlinks.Add($ws.Cells.Item(8,1), "mailto:test3456@styletag.com", "", "", "test3456@styletag.c")
$ws.Cells.Item(8,1).Font.Underline = $false
$ws.Cells.Item(8,1).Font.Color = 16711680
$ws.Cells.Item(8,1).Font.Name = "Arial"

$ws.Hyperlinks.Add($ws.Cells.Item(10,1), "mailto:test3456@styletag.com", "", "", "test3456@styletag.com")
$ws.Cells.Item(10,1).Font.Underline = $false
$ws.Cells.Item(10,1).Font.Color = 16711680
$ws.Cells.Item(10,1).Font.Name = "Arial"

$ws.Hyperlinks.Add($ws.Cells.Item(14,1), "mailto:test3456@styletag.com", "", "", "test3456@styletag.com")
$ws.Cells.Item(14,1).Font.Underline = $false
$ws.Cells.Item(14,1).Font.Color = 16711680
$ws.Cells.Item(14,1).Font.Name = "Arial"

$ws.Hyperlinks.Add($ws.Cells.Item(15,1), "mailto:test34567@styletag.com", "", "", "test34567@styletag.com")
$ws.Cells.Item(15,1).Font.Underline = $false
$ws.Cells.Item(15,1).Font.Color = 16711680
$ws.Cells.Item(15,1).Font.Name = "Arial"

$ws.Hyperlinks.Add($ws.Cells.Item(17,1), "mailto:test3456@styletag.com", "", "", "test3456@styletag.com")
$ws.Cells.Item(17,1).Font.Underline = $false
$ws.Cells.Item(17,1).Font.Color = 16711680
$ws.Cells.Item(17,1).Font.Name = "Arial"
